#
# Locate the run "REGISTER NO: 312214525" and split the value run into
# two runs: ": " and "312214525(80B2F9C1296C9EDB1C547A282A881FA6)".
#
# This mirrors the author's edit (REGISTER NO value annotated with a
# hash suffix) while also reproducing the run-split observed in the
# canonical OOXML diff (the original single run is broken into a ": "
# run and a "312214525(...)" run).

$p = $ppt.ActivePresentation

$registerNumber = "312214525"
$suffix = "(80B2F9C1296C9EDB1C547A282A881FA6)"

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -like "*$registerNumber*") {
                    $targetShape = $shape
                    $targetSlide = $slide
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange

    # Find the register-number run so we don't depend on hard-coded offsets.
    $numberRange = $tr.Find($registerNumber, 0)

    # The ": " immediately precedes the number in the same run.
    $colonStart = $numberRange.Start - 2
    $colonRange = $tr.Characters($colonStart, 2)
    $colonRange.Text = ": "

    # Re-locate the number (the preceding split may shift nothing here,
    # but re-querying keeps this robust) and append the hash suffix.
    $numberRange = $tr.Find($registerNumber, 0)
    $numberOnly = $tr.Characters($numberRange.Start, $numberRange.Length)
    $numberOnly.Text = $registerNumber + $suffix
}
